# Updated symbol list on Mon Dec 19 11:19:33 UTC 2022 with GitHub Actions
# Applies the price/ranking refresh described by the diff:
#  - Column D ("Price") values refreshed for many rows.
#  - Rows 14-25 ("Coin"/"Link"/"Price"/"Volume(1h)") reshuffled to reflect
#    the new coin ranking order (MCDex/BitForexToken swap places; One drops
#    from rank 16 down to rank 25, shifting TigerCash..BTSEToken up by one).
#
# Note: D-column values are numeric-looking strings that must remain text
# (as in the source file, which stores every Price cell as inline/shared
# text). A leading apostrophe forces Excel to keep them as text instead of
# auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price-only updates -------------------------------------------
$ws.Range("D2").Value  = "'248.18"
$ws.Range("D3").Value  = "'21.61"
$ws.Range("D4").Value  = "'5.506"
$ws.Range("D6").Value  = "'3.387"
$ws.Range("D7").Value  = "'0.8072"
$ws.Range("D8").Value  = "'1.038"
$ws.Range("D9").Value  = "'0.1500"
$ws.Range("D10").Value = "'0.07984"
$ws.Range("D11").Value = "'0.03149"
$ws.Range("D13").Value = "'0.09294"

# --- Row 14: was BitForexToken, now MCDex ---------------------------------
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'3.501"
$ws.Range("E14").Value = "13MCDexMCB"

# --- Row 15: was MCDex, now BitForexToken ---------------------------------
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001645"
$ws.Range("E15").Value = "14BitForexTokenBF"

# --- Row 16: CoinExToken stays, only price updates ------------------------
$ws.Range("D16").Value = "'0.04707"

# --- Row 17: was One, now TigerCash ---------------------------------------
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006353"
$ws.Range("E17").Value = "16TigerCashTCH"

# --- Row 18: was TigerCash, now HotbitToken -------------------------------
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").Value = "'0.005042"
$ws.Range("E18").Value = "17HotbitTokenHTB"

# --- Row 19: was HotbitToken, now BitKan ----------------------------------
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Value = "'0.001044"
$ws.Range("E19").Value = "18BitKanKAN"

# --- Row 20: was BitKan, now NitroEx --------------------------------------
$ws.Range("B20").Value = "NitroEx"
$ws.Range("C20").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D20").Value = "'0.0001501"
$ws.Range("E20").Value = "19NitroExNTX"

# --- Row 21: was NitroEx, now UpBots ---------------------------------------
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").Value = "'0.0003202"
$ws.Range("E21").Value = "20UpBotsUBXT"

# --- Row 22: was UpBots, now LEO --------------------------------------------
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.768"
$ws.Range("E22").Value = "21LEOLEO"

# --- Row 23: was LEO, now KuCoinToken ---------------------------------------
$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D23").Value = "'6.425"
$ws.Range("E23").Value = "22KuCoinTokenKCS"

# --- Row 24: was KuCoinToken, now BTSEToken ---------------------------------
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.129"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- Row 25: was BTSEToken, now One -----------------------------------------
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01158"
$ws.Range("E25").Value = "24OneONEBestin24h"

# --- More simple price-only updates ----------------------------------------
$ws.Range("D26").Value = "'0.3279"
$ws.Range("D27").Value = "'0.1301"
$ws.Range("D40").Value = "'0.04109"
$ws.Range("D41").Value = "'0.006973"
$ws.Range("D42").Value = "'0.003502"
$ws.Range("D43").Value = "'0.1043"
$ws.Range("D44").Value = "'0.009141"
$ws.Range("D45").Value = "'0.00005898"
$ws.Range("D47").Value = "'0.0005503"
$ws.Range("D48").Value = "'0.6829"
$ws.Range("D49").Value = "'0.008271"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D51").Value = "'0.01011"
